$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: phone number, entered as text with a quote prefix instead of a plain number
$ws.Range("D2").Value = "'+21626373656"

# E2: color value changes from orange to blue
$ws.Range("E2").Value = "blue"

# Update the active selection to D3
$ws.Range("D3").Select()
